$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.14 = 11730.96 pesos`n✅ 11730.96 pesos = 3.12 = 970.56 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 318.5
$wsTasas.Range("O10").Value = 3736.31
$wsTasas.Range("N12").Value = 3759
$wsTasas.Range("O12").Value = 311
